$d = $word.ActiveDocument

# The document currently ends with two empty paragraphs, right before the
# section break (".../git push -u origin main" paragraph, then two blank
# paragraphs). This edit replaces the very last (trailing) empty paragraph
# with three new paragraphs - numbered "3.", "4." and "5." - describing
# further git branch / checkout workflow steps. The blank paragraph before
# it is left untouched.

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)

if ($lastPara.Range.Text -ne "" -and $lastPara.Range.Text -ne "`r") {
    throw "Expected the trailing paragraph to be empty before inserting the new content."
}

$targetRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$newParagraphsXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">3. check if working tree is clean or </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>not ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> using </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>git status</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">4. Now creating a new branch using </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">git branch </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>[</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>branchname</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>]</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>eg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">git branch </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mybranch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">5.Now checking out the branch made in step </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>4</w:t></w:r><w:r><w:t xml:space="preserve"> ,</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> git checkout </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mybranch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. (Checkout is basically switching to a new branch)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($newParagraphsXml) | Out-Null

# `InsertXML` inserts the new paragraphs just ahead of the original empty
# paragraph mark rather than replacing it outright (Word never lets you
# truly delete the terminal mark of a range this way - it always survives
# as a fresh, still-empty paragraph after the inserted content). Collapse
# that now-redundant mark back into the last newly-inserted paragraph by
# deleting the single character right before it - this merges the two,
# leaving the document with exactly the three new paragraphs where the
# old trailing empty paragraph used to be.
$newLastIndex = $d.Paragraphs.Count
$lastInsertedPara = $d.Paragraphs.Item($newLastIndex - 1)
$mergeRange = $d.Range($lastInsertedPara.Range.End - 1, $lastInsertedPara.Range.End)
$mergeRange.Delete() | Out-Null
